# agenda.pptx — "Before we start" slide (slide 2):
#   * "Any setup issues ?" -> merge the trailing "?" run into "issues?"
#   * "You must have Admin access to " / "your machine" -> merge into one run
#   * "ourse materials" -> split into "ourse " + "materials" runs
#   * swap the Dropbox zip hyperlink for the plain GitHub repo URL (no more link/underline)
#   * "Lab structure" -> split into "Lab " + "structure" runs

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange

# 1) "Any setup issues" + "?" -> "Any setup issues?"
$full = $tr.Text
$joined = "issues?"
$idx = $full.IndexOf($joined)
$run = $tr.Characters($idx + 1, $joined.Length)
$run.Text = $joined

# 2) "You must have Admin access to " + "your machine" -> one sentence
$full = $tr.Text
$sentence = "You must have Admin access to your machine"
$idx = $full.IndexOf($sentence)
$run = $tr.Characters($idx + 1, $sentence.Length)
$run.Text = $sentence

# 3) "ourse materials" -> "ourse " + "materials"
$full = $tr.Text
$piece = "ourse "
$idx = $full.IndexOf($piece)
$run = $tr.Characters($idx + 1, $piece.Length)
$run.Text = $piece

# 4) Point the course-materials link at the GitHub repo instead of the old Dropbox zip,
#    and strip the hyperlink/underline now that it's plain text.
$full = $tr.Text
$oldUrl = "https://www.dropbox.com/s/sbjl67wsgxwaj15/Kraken1.0-node.zip"
$newUrl = "https://github.paypal.com/GlobalTechEd/GlobalTechEd-Node.git"
$idx = $full.IndexOf($oldUrl)
$run = $tr.Characters($idx + 1, $oldUrl.Length)
$run.Text = $newUrl

$run = $tr.Characters($idx + 1, $newUrl.Length)
$run.Font.Underline = 0
$actionSettings = $run.ActionSettings
for ($i = 1; $i -le $actionSettings.Count; $i++) {
    $actionSettings.Item($i).Hyperlink.Address = ""
}

# 5) "Lab structure" -> "Lab " + "structure"
$full = $tr.Text
$piece = "Lab "
$idx = $full.IndexOf("Lab structure")
$run = $tr.Characters($idx + 1, $piece.Length)
$run.Text = $piece
